$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value = 7824.95
$ws.Range("I76").Value = 9928.5
$ws.Range("J76").Value = 2916.6667
$ws.Range("K76").Value = 9928.5
$ws.Range("L76").Value = 2916.6667
$ws.Range("M76").Value = -9613.5
$ws.Range("N76").Value = -3546.6667

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value = 7824.95
$ws.Range("I79").Value = 9928.5
$ws.Range("J79").Value = 2916.6667
$ws.Range("K79").Value = 9928.5
$ws.Range("L79").Value = 2916.6667
$ws.Range("M79").Value = -8836.5
$ws.Range("N79").Value = -5100.6667

# Row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 7780.1055
$ws.Range("I86").Value = 10942.091
$ws.Range("J86").Value = 3432.375
$ws.Range("K86").Value = 10942.091
$ws.Range("L86").Value = 3432.375
$ws.Range("M86").Value = -9819.091
$ws.Range("N86").Value = -5678.375

# Row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 7780.1055
$ws.Range("I89").Value = 10942.091
$ws.Range("J89").Value = 3432.375
$ws.Range("K89").Value = 54710.455
$ws.Range("L89").Value = 17161.875
$ws.Range("M89").Value = -49094.455
$ws.Range("N89").Value = -28393.875

# Row 103: Let Loose the Juice | Persimmon Tannin
$ws.Range("H103").Value = 10526857
$ws.Range("I103").Value = 390
$ws.Range("J103").Value = 22222932
$ws.Range("K103").Value = 1170
$ws.Range("L103").Value = 66668796
$ws.Range("M103").Value = -584
$ws.Range("N103").Value = -66669968

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 1918.2
$ws.Range("I61").Value = 2395.5
$ws.Range("J61").Value = 1600
$ws.Range("K61").Value = 2395.5
$ws.Range("L61").Value = 1600
$ws.Range("M61").Value = -2183.5
$ws.Range("N61").Value = -2024

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1244.9565
$ws.Range("I74").Value = 885.7778
$ws.Range("J74").Value = 2538
$ws.Range("K74").Value = 885.7778
$ws.Range("L74").Value = 2538
$ws.Range("M74").Value = -11.77779999999996
$ws.Range("N74").Value = -4286

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1244.9565
$ws.Range("I77").Value = 885.7778
$ws.Range("J77").Value = 2538
$ws.Range("K77").Value = 4428.889
$ws.Range("L77").Value = 12690
$ws.Range("M77").Value = -60.88900000000012
$ws.Range("N77").Value = -21426

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2172.9424
$ws.Range("I132").Value = 1361.5957
$ws.Range("J132").Value = 9799.6
$ws.Range("K132").Value = 4084.7871
$ws.Range("L132").Value = 29398.8
$ws.Range("M132").Value = -1554.7871
$ws.Range("N132").Value = -34458.8

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1918.2
$ws.Range("I136").Value = 2395.5
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 7186.5
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -4636.5
$ws.Range("N136").Value = -9900

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 3872.1853
$ws.Range("I20").Value = 5938.0713
$ws.Range("J20").Value = 1647.3846
$ws.Range("K20").Value = 5938.0713
$ws.Range("L20").Value = 1647.3846
$ws.Range("M20").Value = -5691.0713
$ws.Range("N20").Value = -2141.3846

# Row 112: Enlistment Highs | Deepgold Sword
$ws.Range("H112").Value = 33425
$ws.Range("J112").Value = 33425
$ws.Range("L112").Value = 33425
$ws.Range("N112").Value = -36379

# Row 118: Cooking for the Future | Titanbronze Culinary Knife
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 1373.1305
$ws.Range("I134").Value = 1360.9
$ws.Range("J134").Value = 1454.6666
$ws.Range("K134").Value = 4082.7
$ws.Range("L134").Value = 4363.9998
$ws.Range("M134").Value = -1547.7
$ws.Range("N134").Value = -9433.9998

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 1635.2
$ws.Range("I58").Value = 1633.9231
$ws.Range("J58").Value = 1638.8889
$ws.Range("K58").Value = 1633.9231
$ws.Range("L58").Value = 1638.8889
$ws.Range("M58").Value = -1430.9231
$ws.Range("N58").Value = -2044.8889

# Row 86: Birch, Please | Birch Lumber
$ws.Range("H86").Value = 4119.8
$ws.Range("I86").Value = 4328.4287
$ws.Range("J86").Value = 3633
$ws.Range("K86").Value = 4328.4287
$ws.Range("L86").Value = 3633
$ws.Range("M86").Value = -3205.4287
$ws.Range("N86").Value = -5879

# Row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws.Range("H89").Value = 4119.8
$ws.Range("I89").Value = 4328.4287
$ws.Range("J89").Value = 3633
$ws.Range("K89").Value = 21642.1435
$ws.Range("L89").Value = 18165
$ws.Range("M89").Value = -16026.1435
$ws.Range("N89").Value = -29397

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 1635.2
$ws.Range("I136").Value = 1633.9231
$ws.Range("J136").Value = 1638.8889
$ws.Range("K136").Value = 4901.7693
$ws.Range("L136").Value = 4916.6667
$ws.Range("M136").Value = -2351.7693
$ws.Range("N136").Value = -10016.6667

$ws = $wb.Worksheets.Item("CUL")
# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 550219.8
$ws.Range("I132").Value = 877918
$ws.Range("J132").Value = 4056.111
$ws.Range("K132").Value = 7901262
$ws.Range("L132").Value = 36504.999
$ws.Range("M132").Value = -7898732
$ws.Range("N132").Value = -41564.999

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 35187.35
$ws.Range("I70").Value = 37074.375
$ws.Range("J70").Value = 4995
$ws.Range("K70").Value = 37074.375
$ws.Range("L70").Value = 4995
$ws.Range("M70").Value = -36804.375
$ws.Range("N70").Value = -5535

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 35187.35
$ws.Range("I73").Value = 37074.375
$ws.Range("J73").Value = 4995
$ws.Range("K73").Value = 37074.375
$ws.Range("L73").Value = 4995
$ws.Range("M73").Value = -36138.375
$ws.Range("N73").Value = -6867

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3187.7441
$ws.Range("I132").Value = 2544.4119
$ws.Range("J132").Value = 3608.3845
$ws.Range("K132").Value = 7633.2357
$ws.Range("L132").Value = 10825.1535
$ws.Range("M132").Value = -5103.2357
$ws.Range("N132").Value = -15885.1535

# Row 134: Guaranteed Gem | Ihuykanite
$ws.Range("H134").Value = 10318.875
$ws.Range("J134").Value = 10318.875
$ws.Range("L134").Value = 30956.625
$ws.Range("N134").Value = -36026.625

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3254.6765
$ws.Range("I132").Value = 3434.7778
$ws.Range("J132").Value = 3052.0625
$ws.Range("K132").Value = 10304.3334
$ws.Range("L132").Value = 9156.1875
$ws.Range("M132").Value = -7774.3334
$ws.Range("N132").Value = -14216.1875

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 2246.5588
$ws.Range("I136").Value = 1691.375
$ws.Range("J136").Value = 3579
$ws.Range("K136").Value = 5074.125
$ws.Range("L136").Value = 10737
$ws.Range("M136").Value = -2524.125
$ws.Range("N136").Value = -15837

# Row 139: Giving Gatherers Their Gear | Gomphotherium Doublet of Gathering
$ws.Range("H139").Value = 64336
$ws.Range("J139").Value = 64336
$ws.Range("L139").Value = 64336
$ws.Range("N139").Value = -74616

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2288.4849
$ws.Range("I132").Value = 2602.1
$ws.Range("J132").Value = 1806
$ws.Range("K132").Value = 7806.299999999999
$ws.Range("L132").Value = 5418
$ws.Range("M132").Value = -5276.299999999999
$ws.Range("N132").Value = -10478

# Row 139: Cruel Climates | Rroneek Serge Trousers of Gathering
$ws.Range("H139").Value = 91712
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 91712
$ws.Range("N139").Value = -101992
$ws.Range("M139").ClearContents()
